$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "56.834.69"
$ws.Range("E2").Value = "  -3.59%  "

Set-TextValue $ws.Range("D3") "2.530.81"
$ws.Range("E3").Value = "  -4.70%  "

$ws.Range("E4").Value = "  +0.14%  "

Set-TextValue $ws.Range("D5") "511.67"
$ws.Range("E5").Value = "  -2.22%  "

Set-TextValue $ws.Range("D6") "139.14"
$ws.Range("E6").Value = "  -3.62%  "

$ws.Range("E7").Value = "  +0.29%  "

Set-TextValue $ws.Range("D8") "0.556"
$ws.Range("E8").Value = "  -2.84%  "

Set-TextValue $ws.Range("D9") "6.45"
$ws.Range("E9").Value = "  -8.19%  "

Set-TextValue $ws.Range("D10") "0.0989"
$ws.Range("E10").Value = "  -3.52%  "

Set-TextValue $ws.Range("D11") "0.323"
$ws.Range("E11").Value = "  -3.44%  "

$ws.Range("E12").Value = "  +0.01%  "

Set-TextValue $ws.Range("D13") "2.979.14"
$ws.Range("E13").Value = "  -4.51%  "

Set-TextValue $ws.Range("D14") "56.875.31"
$ws.Range("E14").Value = "  -3.51%  "

Set-TextValue $ws.Range("D15") "19.94"
$ws.Range("E15").Value = "  -5.31%  "

$ws.Range("E16").Value = "  -3.33%  "

Set-TextValue $ws.Range("D17") "2.562.17"
$ws.Range("E17").Value = "  -3.68%  "

Set-TextValue $ws.Range("D18") "330.37"
$ws.Range("E18").Value = "  -2.45%  "

Set-TextValue $ws.Range("D19") "4.26"
$ws.Range("E19").Value = "  -2.47%  "

Set-TextValue $ws.Range("D20") "10.03"
$ws.Range("E20").Value = "  -3.27%  "

Set-TextValue $ws.Range("D21") "6.09"
$ws.Range("E21").Value = "  -3.97%  "

$ws.Range("E22").Value = "  +0.10%  "

Set-TextValue $ws.Range("D23") "64.01"
$ws.Range("E23").Value = "  +0.37%  "

Set-TextValue $ws.Range("D24") "0.164"
$ws.Range("E24").Value = "  -0.88%  "

$ws.Range("E25").Value = "  +0.19%  "

Set-TextValue $ws.Range("D26") "0.399"
$ws.Range("E26").Value = "  -4.42%  "

Set-TextValue $ws.Range("D27") "2.658.14"
$ws.Range("E27").Value = "  -4.10%  "

Set-TextValue $ws.Range("D28") "6.89"
$ws.Range("E28").Value = "  -2.40%  "

Set-TextValue $ws.Range("D29") "0.0₃0745"
$ws.Range("E29").Value = "  -7.00%  "

$ws.Range("E30").Value = "  +0.11%  "

Set-TextValue $ws.Range("D31") "6.26"
$ws.Range("E31").Value = "  -6.03%  "

$ws.Range("E32").Value = "  -2.96%  "

Set-TextValue $ws.Range("D33") "148.52"
$ws.Range("E33").Value = "  -0.84%  "

Set-TextValue $ws.Range("D34") "18.40"
$ws.Range("E34").Value = "  -2.32%  "

Set-TextValue $ws.Range("D35") "3.94"
$ws.Range("E35").Value = "  -4.72%  "

$ws.Range("E36").Value = "  -5.68%  "

Set-TextValue $ws.Range("D37") "0.838"
$ws.Range("E37").Value = "  -5.80%  "

Set-TextValue $ws.Range("D38") "35.64"
$ws.Range("E38").Value = "  -3.13%  "

Set-TextValue $ws.Range("D39") "0.818"
$ws.Range("E39").Value = "  -5.80%  "

$ws.Range("E40").Value = "  -4.38%  "

$ws.Range("E41").Value = "  +0.09%  "

Set-TextValue $ws.Range("D42") "3.45"
$ws.Range("E42").Value = "  -3.93%  "

$ws.Range("E43").Value = "  -1.91%  "

Set-TextValue $ws.Range("D44") "10.59"
$ws.Range("E44").Value = "  -0.65%  "

Set-TextValue $ws.Range("D45") "0.572"
$ws.Range("E45").Value = "  -7.30%  "

Set-TextValue $ws.Range("D46") "0.0519"
$ws.Range("E46").Value = "  -2.21%  "

Set-TextValue $ws.Range("D47") "256.70"
$ws.Range("E47").Value = "  -6.81%  "

Set-TextValue $ws.Range("D48") "18.41"
$ws.Range("E48").Value = "  -7.45%  "

Set-TextValue $ws.Range("D49") "1.965.20"
$ws.Range("E49").Value = "  -3.97%  "

# Row 50/51: RenderToken and VeChain swap positions, with RenderToken getting updated price/volume
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D50") "4.54"
$ws.Range("E50").Value = "  -3.61%  "

$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D51") "0.0221"
$ws.Range("E51").Value = "  -3.09%  "